$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.140.08'
$ws.Range("E2").Value = '  -0.47%  '
$ws.Range("D3").Value = '3.455.13'
$ws.Range("E3").Value = '  -1.41%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.33'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.63'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.55%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.84'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.54%  '
$ws.Range("E10").Value = '  -2.23%  '
$ws.Range("E11").Value = '  +1.92%  '
$ws.Range("D12").Value = '4.045.94'
$ws.Range("E12").Value = '  -1.46%  '
$ws.Range("E13").Value = '  +2.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.57'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.38%  '
$ws.Range("D15").Value = '3.472.81'
$ws.Range("E15").Value = '  -1.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000172'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.88%  '
$ws.Range("D17").Value = '63.101.91'
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.44'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.44'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.96%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.15'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '386.54'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.561'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '74.46'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.38%  '
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").Value = '3.587.70'
$ws.Range("E25").Value = '  -1.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000115'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.21%  '
$ws.Range("E27").Value = '  -2.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.69'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.05'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.70%  '
$ws.Range("E31").Value = '  -2.81%  '
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.35'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.33'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.83%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.37'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.14%  '
$ws.Range("E36").Value = '  +2.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.05'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '31.91'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '170.04'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.16%  '
$ws.Range("D40").Value = '3.491.43'
$ws.Range("E40").Value = '  -1.45%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0767'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.81%  '
$ws.Range("E42").Value = '  -1.73%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.48'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.26%  '
$ws.Range("E44").Value = '  -1.63%  '
$ws.Range("E45").Value = '  -3.32%  '
$ws.Range("E46").Value = '  -1.83%  '
$ws.Range("D47").Value = '2.582.53'
$ws.Range("E47").Value = '  -1.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.28'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.87'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.42%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.60'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.19%  '
$ws.Range("E51").Value = '  -0.03%  '
